$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - match formatting of existing header cell H1 (bold,
# centered, bordered) by copying H1's format onto the new header cells.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-22
$values = @{
    2  = @(1, 6)
    3  = @(1, 6)
    4  = @(1, 6)
    5  = @(1, 7)
    6  = @(1, 7)
    7  = @(1, 5)
    8  = @(1, 5)
    9  = @(1, 3)
    10 = @(1, 6)
    11 = @(1, 5)
    12 = @(1, 6)
    13 = @(1, 5)
    14 = @(1, 6)
    15 = @(6, 7)
    16 = @(5, 9)
    17 = @(1, 2)
    18 = @(1, 6)
    19 = @(1, 6)
    20 = @(1, 4)
    21 = @(2, 6)
    22 = @(1, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
